# NIT-9002549750 "Estado de Cuenta" update
# - The "Periodo Mora" list (shared strings) was re-sorted from descending
#   to ascending chronological order, which changes what each row in
#   column E displays (rows 16-32).
# - The outlier "Valor Mora" amount (27083) that was on the first period
#   row moved down to the last period row; every other row keeps 31249.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New ascending order of period labels for rows 16..32 (column E)
$periods = @(
    "1912",
    "2001", "2002", "2003", "2004", "2005", "2006", "2007",
    "2008", "2009", "2010", "2011", "2012",
    "2101", "2102", "2103", "2104"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# Swap the "Valor Mora" value between the first and last data rows
$ws.Range("F16").Value = 31249
$ws.Range("F32").Value = 27083
